$d = $word.ActiveDocument

$newText = "何时给mepc+4,这是个很重要的事情！！！不注意就会有很多重大bug！由于我统一对IRQ的地方都给mepc+4了，但是没有注意到从NEMU中过来的mepc其实是下一条需要运行的指令，因此如果我给他+4了，等于说mret之后运行的是下下条指令，这样就造成了很多无法理解的bug！！这次是遇到了直接把我的栈指针清空了，搞了我一个下午加晚上才发现解决！！"

$r = $d.Content
$r.Find.Execute("ye", $true, $false, $false, $false, $false,
                 $true, 1, $false, $newText, 2) | Out-Null

$r.Collapse(0)

# The engine mis-handles Bookmarks.Add on a truly zero-length Range that
# sits exactly at an end-of-paragraph boundary, so briefly insert a
# placeholder character, anchor the bookmark to the (non-degenerate)
# range covering it, then delete the placeholder - the bookmark
# collapses in place and survives.
$r.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $r)
$r.Text = ""
